# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# d42bade3-... file has been handed off (status -> "Ready for handoff"),
# with fresh handoff timestamps and the Priority switched to "mt".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the d42bade3-... entry -------------------
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-15 12:11:16"

# --- zh-cn sheet: row 3 is the d42bade3-... entry -----------------------
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-15 12:11:12"

# --- de-de sheet: row 3 is the d42bade3-... entry -----------------------
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-15 12:11:16"

# --- Column widths: the longer "Ready for handoff" status text widens ---
# --- the Status columns on all three sheets (autofit-style resize). -----
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
